$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / Row 7: add PRODUCT ("laptop") like row 5 already has ---
$ws.Range("E6").Value = "laptop"
$ws.Range("E7").Value = "laptop"

# --- Row 8: D8 should carry the same "password" styling as D3/D2 (theme10 colored, no underline) ---
$ws.Range("D3").Copy()
$ws.Range("D8").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 10: new test case, browser only (mirrors rows 6/7/9 pattern) ---
$ws.Range("A10").Value = "TC_009_GO_TO_HELP_PAGE"
$ws.Range("B10").Value = "Chrome"

# --- Row 11: new test case with password + email/hyperlink columns (mirrors row 3 / row 8) ---
$ws.Range("A11").Value = "TC_010_Go_TO_CONTACT_US_PAGE"
$ws.Range("B11").Value = "Chrome"

$ws.Range("D11").Value = "Prakruthi123!"
$ws.Range("D3").Copy()
$ws.Range("D11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Hyperlinks.Add($ws.Range("F11"), "mailto:prakruthi.koteshwar@gmail.com", "", "", "prakruthi.koteshwar@gmail.com")
$ws.Range("F3").Copy()
$ws.Range("F11").PasteSpecial(-4122)  # xlPasteFormats, restores the shared Hyperlink style index

# --- Row 12: new test case, browser only ---
$ws.Range("A12").Value = "TC_011_CHANGE_COUNTRY"
$ws.Range("B12").Value = "Chrome"

# --- Row 13: new test case with password + email/hyperlink columns ---
$ws.Range("A13").Value = "TC_012_ADD_NEW_ADDRESS_PAGE"
$ws.Range("B13").Value = "Chrome"

$ws.Range("D13").Value = "Prakruthi123!"
$ws.Range("D3").Copy()
$ws.Range("D13").PasteSpecial(-4122)  # xlPasteFormats

$ws.Hyperlinks.Add($ws.Range("F13"), "mailto:prakruthi.koteshwar@gmail.com", "", "", "prakruthi.koteshwar@gmail.com")
$ws.Range("F3").Copy()
$ws.Range("F13").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# --- The sheet shrank by one empty trailing row ---
$ws.Rows(18).Delete()

# --- Match the saved selection in the workbook ---
$ws.Range("E7").Select()
